$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two data rows that were dropped (RM 232 and SC 92).
# Delete bottom-most row first so the other row index stays valid.
$ws.Rows("28").Delete()
$ws.Rows("26").Delete()

# Apply the individual cell value changes (row numbers below are the
# final row numbers after the two rows above were removed).
$ws.Range("E5").ClearContents()
$ws.Range("F6").Value = 16.43
$ws.Range("E8").Value = -6.6
$ws.Range("F11").Value = 17.65
$ws.Range("E12").ClearContents()
$ws.Range("F12").ClearContents()
$ws.Range("E14").Value = -5.4
$ws.Range("F17").ClearContents()
$ws.Range("E18").ClearContents()
$ws.Range("F25").Value = 16.6
$ws.Range("B26").Value = -20.2
$ws.Range("B27").ClearContents()
$ws.Range("F31").ClearContents()
$ws.Range("F32").ClearContents()
$ws.Range("C33").Value = 10.4
